# "Adding the changes we made on may 9th"
#
# The accelerometer sample table on Sheet1 (header "x","y","z" in row 1,
# data triples in rows 2..21) gains 10 more samples:
#   - 6 new samples belong BEFORE the existing series -> they become the
#     new rows 2-7, and the old rows 2-21 shift down to rows 8-27.
#   - 4 new samples belong AFTER the existing series -> they become the
#     new rows 28-31.
# Net effect: the used range grows from A1:C21 to A1:C31, and every row
# from the original table reappears unchanged, just 6 rows lower.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the 20 rows that were already there (old rows 2-21), so we
# can re-place them 6 rows further down without retyping every number.
$existing = @()
for ($r = 2; $r -le 21; $r++) {
    $existing += @(, @($ws.Cells.Item($r, 1).Value(), $ws.Cells.Item($r, 2).Value(), $ws.Cells.Item($r, 3).Value()))
}

# The 6 brand-new leading samples (new rows 2-7).
$leading = @(
    @(-3.195676267147064, 5.127160429954529, -1.443197593092918),
    @(-3.130342268943787, 5.136516356468201, -1.369547128677368),
    @(-3.034864258766174, 5.101877164840698, -1.325036150217056),
    @(-3.194309616088868, 5.024436473846436, -1.315180826187134),
    @(-3.382834231853486, 5.098868799209595, -1.453447324037552),
    @(-3.195986032485962, 5.139615774154663, -1.564420849084854)
)

# The 4 brand-new trailing samples (new rows 28-31).
$trailing = @(
    @(2.14622653722763,  5.513726615905762, -1.284043130278588),
    @(2.037818813323974, 5.183717918395995, -1.269947481155395),
    @(2.045576536655426, 5.118093979358673, -1.376126399636268),
    @(2.174056196212769, 5.255697178840638, -1.429987555742264)
)

$allRows = @()
$allRows += $leading
$allRows += $existing
$allRows += $trailing

# Write the full, re-numbered series back starting at row 2. Writing the
# cell values directly (rather than calling Rows.Insert) keeps every cell
# on its original, unstyled number format - exactly like the rows that
# were already in the sheet.
for ($i = 0; $i -lt $allRows.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $allRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $allRows[$i][1]
    $ws.Cells.Item($row, 3).Value = $allRows[$i][2]
}

Write-Output ("Updated range: " + $ws.UsedRange.Address())
